$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.140.71"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "3.330.21"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.17"
$ws.Range("E5").Value = "  +3.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.64"
$ws.Range("E6").Value = "  -2.45%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.325.29"
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("E9").Value = "  -2.47%  "
$ws.Range("E10").Value = "  -2.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.580"
$ws.Range("E11").Value = "  -1.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.13"
$ws.Range("E12").Value = "  -1.78%  "
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "677.55"
$ws.Range("E14").Value = "  +10.63%  "
$ws.Range("D15").Value = "3.857.07"
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("E16").Value = "  -2.65%  "
$ws.Range("D17").Value = "66.238.81"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.91"
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.118"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").Value = "3.329.28"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.13"
$ws.Range("E21").Value = "  -0.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.898"
$ws.Range("E22").Value = "  -2.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.92"
$ws.Range("E23").Value = "  -3.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.35"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.04"
$ws.Range("E25").Value = "  -2.57%  "
$ws.Range("E26").Value = "  -1.40%  "
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("E28").Value = "  -2.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.26"
$ws.Range("E29").Value = "  +5.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.50"
$ws.Range("E30").Value = "  -2.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.78"
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "602.90"
$ws.Range("E32").Value = "  +5.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.90"
$ws.Range("E33").Value = "  -4.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.00"
$ws.Range("E34").Value = "  -1.51%  "
$ws.Range("E35").Value = "  -1.07%  "
$ws.Range("D36").Value = "3.814.74"
$ws.Range("E36").Value = "  +1.68%  "
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "56.09"
$ws.Range("E38").Value = "  -2.49%  "
$ws.Range("E39").Value = "  -2.10%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0699"
$ws.Range("E40").Value = "  -4.67%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.127"
$ws.Range("E41").Value = "  -3.79%  "
$ws.Range("E42").Value = "  -4.12%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.19"
$ws.Range("E43").Value = "  -4.31%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.41"
$ws.Range("E44").Value = "  +5.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.337"
$ws.Range("E45").Value = "  -1.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0414"
$ws.Range("E46").Value = "  -2.69%  "
$ws.Range("E47").Value = "  -11.79%  "
$ws.Range("E48").Value = "  -2.08%  "
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.56"
$ws.Range("E50").Value = "  -2.76%  "
$ws.Range("E51").Value = "  +1.14%  "
